$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.492.50"
$ws.Range("E2").Value = "  +4.59%  "

$ws.Range("D3").Value = "3.129.46"
$ws.Range("E3").Value = "  +3.39%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.50"
$ws.Range("E5").Value = "  +7.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.72"
$ws.Range("E6").Value = "  +2.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.376"
$ws.Range("E7").Value = "  +5.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.967"
$ws.Range("E8").Value = "  +20.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "3.125.42"
$ws.Range("E10").Value = "  +3.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.753"
$ws.Range("E11").Value = "  +29.47%  "

$ws.Range("E12").Value = "  +7.10%  "

$ws.Range("E13").Value = "  +9.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.90"
$ws.Range("E14").Value = "  +12.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.46"
$ws.Range("E15").Value = "  +5.33%  "

$ws.Range("D16").Value = "91.173.14"
$ws.Range("E16").Value = "  +4.50%  "

$ws.Range("D17").Value = "3.705.92"
$ws.Range("E17").Value = "  +3.32%  "

$ws.Range("D18").Value = "3.132.56"
$ws.Range("E18").Value = "  +4.12%  "

$ws.Range("E19").Value = "  +20.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000226"
$ws.Range("E20").Value = "  +18.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.17"
$ws.Range("E21").Value = "  +10.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "435.27"
$ws.Range("E22").Value = "  +5.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.76"
$ws.Range("E23").Value = "  +11.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.18"
$ws.Range("E24").Value = "  +8.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.85"
$ws.Range("E25").Value = "  +12.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.43"
$ws.Range("E26").Value = "  +8.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "80.33"
$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("D28").Value = "3.295.40"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.170"
$ws.Range("E30").Value = "  +8.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.00"
$ws.Range("E31").Value = "  +14.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "534.00"
$ws.Range("E32").Value = "  +8.25%  "

$ws.Range("E33").Value = "  +14.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.887"
$ws.Range("E34").Value = "  -18.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.27"
$ws.Range("E35").Value = "  +13.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.145"
$ws.Range("E36").Value = "  +12.06%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "23.73"
$ws.Range("E37").Value = "  +9.34%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.30"
$ws.Range("E38").Value = "  +8.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.86"
$ws.Range("E39").Value = "  +5.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.32"

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.154"
$ws.Range("E42").Value = "  +17.78%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0778"
$ws.Range("E43").Value = "  +20.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.381"
$ws.Range("E44").Value = "  +8.37%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("E46").Value = "  +10.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "145.00"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.23"
$ws.Range("E48").Value = "  +2.64%  "

$ws.Range("E49").Value = "  +13.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000263"
$ws.Range("E50").Value = "  +28.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "167.38"
$ws.Range("E51").Value = "  +10.16%  "
